$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 570-571, shifting the existing data (old rows 570-676)
# down to become rows 572-678.
$ws.Rows("570:571").Insert()

# New row 570: Coliflor, Primera, week of 2022-03-24 (serial 44644)
$ws.Cells(570,1).Value = 6
$ws.Cells(570,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells(570,3).Value = "Metropolitana"
$ws.Cells(570,4).Value = 44644
$ws.Cells(570,5).Value = 13
$ws.Cells(570,6).Value = 100112008
$ws.Cells(570,7).Value = "Coliflor"
$ws.Cells(570,8).Value = "Sin especificar"
$ws.Cells(570,9).Value = "Primera"
$ws.Cells(570,10).Value = 8700
$ws.Cells(570,11).Value = 1000
$ws.Cells(570,12).Value = 1200
$ws.Cells(570,13).Value = 1124
$ws.Cells(570,14).Value = "`$/unidad"
$ws.Cells(570,15).Value = "Región Metropolitana"
$ws.Cells(570,16).Value = 1124
$ws.Cells(570,17).Value = 1
$ws.Cells(570,18).Value = "Hortaliza"

# New row 571: Coliflor, Segunda, same week (serial 44644)
$ws.Cells(571,1).Value = 6
$ws.Cells(571,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells(571,3).Value = "Metropolitana"
$ws.Cells(571,4).Value = 44644
$ws.Cells(571,5).Value = 13
$ws.Cells(571,6).Value = 100112008
$ws.Cells(571,7).Value = "Coliflor"
$ws.Cells(571,8).Value = "Sin especificar"
$ws.Cells(571,9).Value = "Segunda"
$ws.Cells(571,10).Value = 3800
$ws.Cells(571,11).Value = 700
$ws.Cells(571,12).Value = 800
$ws.Cells(571,13).Value = 734
$ws.Cells(571,14).Value = "`$/unidad"
$ws.Cells(571,15).Value = "Región Metropolitana"
$ws.Cells(571,16).Value = 734
$ws.Cells(571,17).Value = 1
$ws.Cells(571,18).Value = "Hortaliza"

# Apply the same date number format used by the other D-column (Fecha) cells
# to the two new D cells, matching the surrounding rows.
$ws.Cells(570,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(571,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
